$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (k value)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary labels (column A) and aggregate formulas (column B), rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold/size-12/vertical-centered style on a scratch cell, then
# broadcast it to B14:B17 via copy/paste-special so only a single new
# font + cellXf gets minted (avoids orphaned intermediate style records).
$tmpl = $ws.Range("D100")
$tmpl.Font.Bold = $true
$tmpl.Font.Size = 12
$tmpl.VerticalAlignment = -4108
$tmpl.Copy()
$ws.Range("B14:B17").PasteSpecial(-4122)
$tmpl.Clear()

$ws.Range("A14:B17").RowHeight = 15.6

$ws.Range("A14:B17").Select() | Out-Null

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
